# Se actualiza los datos de pagos
# Updates the payment test-data values that live in the shared-string
# table and are referenced from the "numeroUsuario", "placa" and
# "vigencia" columns of the data-driven sheets, plus the remembered
# selection on the first sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "SucripcionDesdePagos" (1st sheet) ---------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "'72934725"
$ws1.Range("F2").Value = "'XFN-363"
$ws1.Range("G2").Value = "'06/11/2025"

# --- Sheet "SucripcionDesdeAfiliacion" (3rd sheet) -----------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "'72934725"
$ws3.Range("F2").Value = "'XFN-363"
$ws3.Range("G2").Value = "'06/11/2025"

# --- Sheet "SucripcionDesdeOpcionPagar" (4th sheet) -----------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = "'72934725"
$ws4.Range("C2").Value = "'XFN-363"
$ws4.Range("G2").Value = "'06/11/2025"

# --- Remembered selection on sheet 1 moves from G2 to F2:G2 -------------
[void]$ws1.Range("F2:G2").Select()

# Restore the originally active sheet/tab (4th sheet) so the workbook's
# active-tab/tabSelected bookkeeping is left untouched.
[void]$ws4.Activate()
